$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Section/course label
$ws.Range("B3").Value = "BSCS 3 IRREG"

# Row 6
$ws.Range("A6").Value = "NSCI"
$ws.Range("B6").Value = 6100
$ws.Range("C6").Value = "Calculus-based Physics 1"
$ws.Range("D6").Value = 3
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 4
$ws.Range("H6").Value = 6088
$ws.Range("I6").Value = 6089

# Row 7
$ws.Range("A7").Value = "FILI"
$ws.Range("B7").Value = 6201
$ws.Range("C7").Value = "Kritikal na Pagbasa Pagsulat at Pagsasalita"
$ws.Range("D7").Value = 3
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 3
$ws.Range("H7").Value = 1666

# Row 8
$ws.Range("A8").Value = "CS"
$ws.Range("B8").Value = 6300
$ws.Range("C8").Value = "Software Engineering 2"
$ws.Range("D8").Value = 2
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 3
$ws.Range("H8").Value = 1671
$ws.Range("I8").Value = 1672

# Row 9 (taller row)
$ws.Range("A9").Value = "CS"
$ws.Range("B9").Value = 6206
$ws.Range("C9").Value = "Principles of Operating Systems and its Application"
$ws.Range("D9").Value = 2
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 3
$ws.Range("H9").Value = 1715
$ws.Range("I9").Value = 1716
$ws.Rows.Item(9).RowHeight = 31.5

# Row 10
$ws.Range("A10").Value = "ITE"
$ws.Range("B10").Value = 6301
$ws.Range("C10").Value = "Technopreneurship"
$ws.Range("D10").Value = 3
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 3
$ws.Range("H10").Value = 1717

# Row 11
$ws.Range("A11").Value = "GE"
$ws.Range("B11").Value = 6101
$ws.Range("C11").Value = "Readings in Philippine History"
$ws.Range("D11").Value = 3
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 3
$ws.Range("H11").Value = 1678

# Row 12
$ws.Range("A12").Value = "IT"
$ws.Range("B12").Value = 6316
$ws.Range("C12").Value = "CS Major Elective 3-XML-Based Web Applications"
$ws.Range("D12").Value = 2
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 3
$ws.Range("H12").Value = 1680
$ws.Range("I12").Value = 1681

# Row 13
$ws.Range("A13").Value = "CS"
$ws.Range("B13").Value = 6398
$ws.Range("C13").Value = "CS Design Project 1"
$ws.Range("D13").Value = 3
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 3
$ws.Range("H13").Value = 1682

# Row 14 (was entirely blank)
$ws.Range("A14").Value = "GE"
$ws.Range("B14").Value = 6301
$ws.Range("C14").Value = "Life and Work of Rizal"
$ws.Range("D14").Value = 3
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 3
$ws.Range("H14").Value = 1684

# Selection as last-seen in file
$ws.Range("E16").Select() | Out-Null
